$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2024-04-21 Sunday" "2024-04-22 Monday"

Replace-Text "85×68=" "77×76="
Replace-Text "36×50=" "24×17="
Replace-Text "64×20=" "49×87="
Replace-Text "14×79=" "37×27="
Replace-Text "44×79=" "71×55="
Replace-Text "88×82=" "35×85="
Replace-Text "71×69=" "44×99="
Replace-Text "85×28=" "28×85="
Replace-Text "96×72=" "98×34="
Replace-Text "79×95=" "54×84="
Replace-Text "47×11=" "32×81="
Replace-Text "52×93=" "42×72="
Replace-Text "80×89=" "72×40="
Replace-Text "76×56=" "54×69="
Replace-Text "24×19=" "66×25="
Replace-Text "42×20=" "63×37="
Replace-Text "35×94=" "42×76="
Replace-Text "41×92=" "82×64="
Replace-Text "35×55=" "73×32="
Replace-Text "72×35=" "79×42="
Replace-Text "74×14=" "52×68="
Replace-Text "48×50=" "14×67="
Replace-Text "98×47=" "40×81="
Replace-Text "77×64=" "55×68="
Replace-Text "34×21=" "82×60="

Write-Output "Done"
